$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Gdf1"
$ws.Range("C2").Value = "Bmpr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.00687
$ws.Range("H2").Value = 0.02061
$ws.Range("I2").Value = 0.2011300758263314
$ws.Range("J2").Value = 0.2011300758263313
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 30.46625333333334
$ws.Range("N2").Value = 91.39876000000001
$ws.Range("O2").Value = 0.2185380492512374
$ws.Range("P2").Value = 0.2331534018544084
$ws.Range("Q2").Value = 0.2093031604
$ws.Range("R2").Value = 1.8837284436
$ws.Range("S2").Value = 0.04395457441683992
$ws.Range("T2").Value = 0.04689416139414427

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gdf1"
$ws.Range("C3").Value = "Bmpr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.00687
$ws.Range("H3").Value = 0.02061
$ws.Range("I3").Value = 0.2011300758263314
$ws.Range("J3").Value = 0.2011300758263313
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 34.738136
$ws.Range("N3").Value = 104.214408
$ws.Range("O3").Value = 0.2491807703757967
$ws.Range("P3").Value = 0.2658454419670822
$ws.Range("Q3").Value = 0.23865099432
$ws.Range("R3").Value = 2.14785894888
$ws.Range("S3").Value = 0.05011774724014765
$ws.Range("T3").Value = 0.05346951390092381

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gdf1"
$ws.Range("C4").Value = "Bmpr2"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.00687
$ws.Range("H4").Value = 0.02061
$ws.Range("I4").Value = 0.2011300758263314
$ws.Range("J4").Value = 0.2011300758263313
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 23.69037333333334
$ws.Range("N4").Value = 71.07112000000001
$ws.Range("O4").Value = 0.1699338582153697
$ws.Range("P4").Value = 0.181298667526812
$ws.Range("Q4").Value = 0.1627528648
$ws.Range("R4").Value = 1.4647757832
$ws.Range("S4").Value = 0.03417880978831836
$ws.Range("T4").Value = 0.03646461474688053

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gdf1"
$ws.Range("C5").Value = "Bmpr2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.00687
$ws.Range("H5").Value = 0.02061
$ws.Range("I5").Value = 0.2011300758263314
$ws.Range("J5").Value = 0.2011300758263313
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.2168665
$ws.Range("N5").Value = 52.433733
$ws.Range("O5").Value = 0.18805669340777
$ws.Range("P5").Value = 0.1337556791894743
$ws.Range("Q5").Value = 0.180109872855
$ws.Range("R5").Value = 1.08065923713
$ws.Range("S5").Value = 0.03782385700475392
$ws.Range("T5").Value = 0.02690228989758141

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gdf1"
$ws.Range("C6").Value = "Bmpr2"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.00687
$ws.Range("H6").Value = 0.02061
$ws.Range("I6").Value = 0.2011300758263314
$ws.Range("J6").Value = 0.2011300758263313
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 24.297748
$ws.Range("N6").Value = 72.893244
$ws.Range("O6").Value = 0.1742906287498262
$ws.Range("P6").Value = 0.1859468094622229
$ws.Range("Q6").Value = 0.16692552876
$ws.Range("R6").Value = 1.50232975884
$ws.Range("S6").Value = 0.03505508737627151
$ws.Range("T6").Value = 0.03739949588680128

# Row 7
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Gdf1"
$ws.Range("C7").Value = "Bmpr2"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.027287
$ws.Range("H7").Value = 0.081861
$ws.Range("I7").Value = 0.7988699241736688
$ws.Range("J7").Value = 0.7988699241736686
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 30.46625333333334
$ws.Range("N7").Value = 91.39876000000001
$ws.Range("O7").Value = 0.2185380492512374
$ws.Range("P7").Value = 0.2331534018544084
$ws.Range("Q7").Value = 0.8313326547066668
$ws.Range("R7").Value = 7.481993892360001
$ws.Range("S7").Value = 0.1745834748343975
$ws.Range("T7").Value = 0.1862592404602642

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Gdf1"
$ws.Range("C8").Value = "Bmpr2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.027287
$ws.Range("H8").Value = 0.081861
$ws.Range("I8").Value = 0.7988699241736688
$ws.Range("J8").Value = 0.7988699241736686
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 34.738136
$ws.Range("N8").Value = 104.214408
$ws.Range("O8").Value = 0.2491807703757967
$ws.Range("P8").Value = 0.2658454419670822
$ws.Range("Q8").Value = 0.947899517032
$ws.Range("R8").Value = 8.531095653288
$ws.Range("S8").Value = 0.1990630231356491
$ws.Range("T8").Value = 0.2123759280661584

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Gdf1"
$ws.Range("C9").Value = "Bmpr2"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.027287
$ws.Range("H9").Value = 0.081861
$ws.Range("I9").Value = 0.7988699241736688
$ws.Range("J9").Value = 0.7988699241736686
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 23.69037333333334
$ws.Range("N9").Value = 71.07112000000001
$ws.Range("O9").Value = 0.1699338582153697
$ws.Range("P9").Value = 0.181298667526812
$ws.Range("Q9").Value = 0.6464392171466669
$ws.Range("R9").Value = 5.817952954320001
$ws.Range("S9").Value = 0.1357550484270514
$ws.Range("T9").Value = 0.1448340527799315

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Gdf1"
$ws.Range("C10").Value = "Bmpr2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.027287
$ws.Range("H10").Value = 0.081861
$ws.Range("I10").Value = 0.7988699241736688
$ws.Range("J10").Value = 0.7988699241736686
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 26.2168665
$ws.Range("N10").Value = 52.433733
$ws.Range("O10").Value = 0.18805669340777
$ws.Range("P10").Value = 0.1337556791894743
$ws.Range("Q10").Value = 0.7153796361855002
$ws.Range("R10").Value = 4.292277817113001
$ws.Range("S10").Value = 0.1502328364030161
$ws.Range("T10").Value = 0.1068533892918929

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Gdf1"
$ws.Range("C11").Value = "Bmpr2"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.027287
$ws.Range("H11").Value = 0.081861
$ws.Range("I11").Value = 0.7988699241736688
$ws.Range("J11").Value = 0.7988699241736686
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 24.297748
$ws.Range("N11").Value = 72.893244
$ws.Range("O11").Value = 0.1742906287498262
$ws.Range("P11").Value = 0.1859468094622229
$ws.Range("Q11").Value = 0.663012649676
$ws.Range("R11").Value = 5.967113847084
$ws.Range("S11").Value = 0.1392355413735547
$ws.Range("T11").Value = 0.1485473135754217
